# fix(publipostage): Correct status name
#
# Renames a handful of "statut_label" / "statut_name" shared-string values
# used throughout the mailmerge source sheet:
#   "bleu"                                                   -> "noir"
#   "pas de résultat ni de publication"                      -> "pas de résultat postés ni publiés"
#   "résultat et / ou publication posté"                     -> "résultat postés ou publiés"
#   "résultat et / ou publication posté dans les 36 mois"    -> "résultat postés ou publiés dans les 36 mois"
#   "résultat et / ou publication posté dans les 12 mois"    -> "résultat postés ou publiés dans les 12 mois"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1
$xlByRows = 1

$map = @{
    "bleu" = "noir"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
}

foreach ($old in $map.Keys) {
    $new = $map[$old]
    $ws.Cells.Replace($old, $new, $xlWhole, $xlByRows, $false, $false, $false, $false)
}
